$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the existing header cell (H1) into the two new header
# cells so they get the same bold/bordered/centered formatting (style index 1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header row text
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-14
$values = @(
    @(2, 8, 9),
    @(3, 9, 9),
    @(4, 9, 9),
    @(5, 9, 9),
    @(6, 4, 6),
    @(7, 7, 8),
    @(8, 7, 7),
    @(9, 7, 7),
    @(10, 8, 8),
    @(11, 9, 9),
    @(12, 1, 2),
    @(13, 6, 7),
    @(14, 3, 3)
)

foreach ($row in $values) {
    $r = $row[0]
    $i = $row[1]
    $j = $row[2]
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
}
